$wb = $excel.ActiveWorkbook

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3951.1667
$ws.Range("I64").Value = 3519.2942
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 3519.2942
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -3271.2942
$ws.Range("N64").Value = -5496

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3951.1667
$ws.Range("I67").Value = 3519.2942
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 3519.2942
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -2661.2942
$ws.Range("N67").Value = -6716

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4338.4614
$ws.Range("I74").Value = 3750
$ws.Range("J74").Value = 4445.4546
$ws.Range("K74").Value = 3750
$ws.Range("L74").Value = 4445.4546
$ws.Range("M74").Value = -2814
$ws.Range("N74").Value = -6317.4546

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4338.4614
$ws.Range("I77").Value = 3750
$ws.Range("J77").Value = 4445.4546
$ws.Range("K77").Value = 18750
$ws.Range("L77").Value = 22227.273
$ws.Range("M77").Value = -14070
$ws.Range("N77").Value = -31587.273

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 74076900
$ws.Range("I100").Value = 27779846
$ws.Range("J100").Value = 166671000
$ws.Range("K100").Value = 27779846
$ws.Range("L100").Value = 166671000
$ws.Range("M100").Value = -27779305
$ws.Range("N100").Value = -166672082

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2094.9048
$ws.Range("I138").Value = 1403.5476
$ws.Range("J138").Value = 3477.6191
$ws.Range("K138").Value = 4210.642800000001
$ws.Range("L138").Value = 10432.8573
$ws.Range("M138").Value = 929.3571999999995
$ws.Range("N138").Value = -20712.8573

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1908.1333
$ws.Range("I2").Value = 2008.6154
$ws.Range("J2").Value = 1255
$ws.Range("K2").Value = 2008.6154
$ws.Range("L2").Value = 1255
$ws.Range("M2").Value = -1895.6154
$ws.Range("N2").Value = -1481

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 945.6923
$ws.Range("I97").Value = 677.5294
$ws.Range("J97").Value = 1452.2222
$ws.Range("K97").Value = 677.5294
$ws.Range("L97").Value = 1452.2222
$ws.Range("M97").Value = -181.5294
$ws.Range("N97").Value = -2444.2222

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1300
$ws.Range("I102").Value = 1300
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1300
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 322

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1908.1333
$ws.Range("I116").Value = 2008.6154
$ws.Range("J116").Value = 1255
$ws.Range("K116").Value = 2008.6154
$ws.Range("L116").Value = 1255
$ws.Range("M116").Value = 285.3846000000001
$ws.Range("N116").Value = -5843

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1468.2703
$ws.Range("I132").Value = 813.4545000000001
$ws.Range("J132").Value = 3363.7896
$ws.Range("K132").Value = 2440.3635
$ws.Range("L132").Value = 10091.3688
$ws.Range("M132").Value = 89.63649999999961
$ws.Range("N132").Value = -15151.3688

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 37944
$ws.Range("I133").Value = 10000
$ws.Range("J133").Value = 47258.668
$ws.Range("K133").Value = 10000
$ws.Range("L133").Value = 47258.668
$ws.Range("M133").Value = -7470
$ws.Range("N133").Value = -52318.668

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1908.1333
$ws.Range("I3").Value = 2008.6154
$ws.Range("J3").Value = 1255
$ws.Range("K3").Value = 2008.6154
$ws.Range("L3").Value = 1255
$ws.Range("M3").Value = -1894.6154
$ws.Range("N3").Value = -1483

# BSM row 70
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 69459
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 69459
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 69459
$ws.Range("N70").Value = -70045

# BSM row 73
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H73").Value = 69459
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 69459
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 69459
$ws.Range("N73").Value = -71487

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1837.7778
$ws.Range("I86").Value = 2139.4546
$ws.Range("J86").Value = 1549.2174
$ws.Range("K86").Value = 2139.4546
$ws.Range("L86").Value = 1549.2174
$ws.Range("M86").Value = -1016.4546
$ws.Range("N86").Value = -3795.2174

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1837.7778
$ws.Range("I89").Value = 2139.4546
$ws.Range("J89").Value = 1549.2174
$ws.Range("K89").Value = 10697.273
$ws.Range("L89").Value = 7746.087
$ws.Range("M89").Value = -5081.273000000001
$ws.Range("N89").Value = -18978.087

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 806.4167
$ws.Range("I94").Value = 539.4
$ws.Range("J94").Value = 997.1429000000001
$ws.Range("K94").Value = 539.4
$ws.Range("L94").Value = 997.1429000000001
$ws.Range("M94").Value = -88.39999999999998
$ws.Range("N94").Value = -1899.1429

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1409.75
$ws.Range("I99").Value = 995.7917
$ws.Range("J99").Value = 2030.6875
$ws.Range("K99").Value = 995.7917
$ws.Range("L99").Value = 2030.6875
$ws.Range("M99").Value = 502.2083
$ws.Range("N99").Value = -5026.6875

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5052594
$ws.Range("I105").Value = 5052594
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5052594
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -5050847

# BSM row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 45000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 45000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -55140

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1371.0204
$ws.Range("I58").Value = 990.69446
$ws.Range("J58").Value = 2424.2307
$ws.Range("K58").Value = 990.69446
$ws.Range("L58").Value = 2424.2307
$ws.Range("M58").Value = -787.69446
$ws.Range("N58").Value = -2830.2307

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5111.778
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 7751.5
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 7751.5
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -8999.5

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5111.778
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 7751.5
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 38757.5
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -44997.5

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 90910616
$ws.Range("I105").Value = 1677
$ws.Range("J105").Value = 1000000000
$ws.Range("K105").Value = 1677
$ws.Range("L105").Value = 1000000000
$ws.Range("M105").Value = 70
$ws.Range("N105").Value = -1000003494

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1371.0204
$ws.Range("I136").Value = 990.69446
$ws.Range("J136").Value = 2424.2307
$ws.Range("K136").Value = 2972.08338
$ws.Range("L136").Value = 7272.6921
$ws.Range("M136").Value = -422.08338
$ws.Range("N136").Value = -12372.6921

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1827.9
$ws.Range("I4").Value = 299
$ws.Range("J4").Value = 1997.7778
$ws.Range("K4").Value = 897
$ws.Range("L4").Value = 5993.3334
$ws.Range("M4").Value = -785
$ws.Range("N4").Value = -6217.3334

# CUL row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1375
$ws.Range("I64").Value = 1200
$ws.Range("J64").Value = 1900
$ws.Range("K64").Value = 3600
$ws.Range("L64").Value = 5700
$ws.Range("M64").Value = -3330
$ws.Range("N64").Value = -6240

# CUL row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 1375
$ws.Range("I67").Value = 1200
$ws.Range("J67").Value = 1900
$ws.Range("K67").Value = 3600
$ws.Range("L67").Value = 5700
$ws.Range("M67").Value = -2664
$ws.Range("N67").Value = -7572

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5053.5386
$ws.Range("I131").Value = 494.875
$ws.Range("J131").Value = 7079.6113
$ws.Range("K131").Value = 1484.625
$ws.Range("L131").Value = 21238.8339
$ws.Range("M131").Value = 3555.375
$ws.Range("N131").Value = -31318.8339

# GSM row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9001417
$ws.Range("I11").Value = 11111111
$ws.Range("J11").Value = 2672333.2
$ws.Range("K11").Value = 11111111
$ws.Range("L11").Value = 2672333.2
$ws.Range("M11").Value = -11110972
$ws.Range("N11").Value = -2672611.2

# GSM row 12
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 6535333.5
$ws.Range("I12").Value = 7000000
$ws.Range("J12").Value = 3515000
$ws.Range("K12").Value = 7000000
$ws.Range("L12").Value = 3515000
$ws.Range("M12").Value = -6999860
$ws.Range("N12").Value = -3515280

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4918.9814
$ws.Range("I70").Value = 4328.5
$ws.Range("J70").Value = 6606.0713
$ws.Range("K70").Value = 4328.5
$ws.Range("L70").Value = 6606.0713
$ws.Range("M70").Value = -4058.5
$ws.Range("N70").Value = -7146.0713

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4918.9814
$ws.Range("I73").Value = 4328.5
$ws.Range("J73").Value = 6606.0713
$ws.Range("K73").Value = 4328.5
$ws.Range("L73").Value = 6606.0713
$ws.Range("M73").Value = -3392.5
$ws.Range("N73").Value = -8478.0713

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3004.5454
$ws.Range("I80").Value = 2240
$ws.Range("J80").Value = 3641.6667
$ws.Range("K80").Value = 2240
$ws.Range("L80").Value = 3641.6667
$ws.Range("M80").Value = -1242
$ws.Range("N80").Value = -5637.6667

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3004.5454
$ws.Range("I83").Value = 2240
$ws.Range("J83").Value = 3641.6667
$ws.Range("K83").Value = 11200
$ws.Range("L83").Value = 18208.3335
$ws.Range("M83").Value = -6208
$ws.Range("N83").Value = -28192.3335

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1426.3636
$ws.Range("I97").Value = 1641.1111
$ws.Range("J97").Value = 460
$ws.Range("K97").Value = 1641.1111
$ws.Range("L97").Value = 460
$ws.Range("M97").Value = -1145.1111
$ws.Range("N97").Value = -1452

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3666.111
$ws.Range("I102").Value = 3141.6667
$ws.Range("J102").Value = 4715
$ws.Range("K102").Value = 3141.6667
$ws.Range("L102").Value = 4715
$ws.Range("M102").Value = -1519.6667
$ws.Range("N102").Value = -7959

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3000
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2459

# WVR row 12
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 30002.334
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 30002.334
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 30002.334
$ws.Range("N12").Value = -30286.334

# WVR row 20
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 29405.5
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 29405.5
$ws.Range("K20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("M20").Value = 29405.5
$ws.Range("N20").Value = -29885.5

# WVR row 43
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 30000
$ws.Range("I43").Value = 10000
$ws.Range("J43").Value = 50000
$ws.Range("K43").Value = 10000
$ws.Range("L43").Value = 50000
$ws.Range("M43").Value = -9851
$ws.Range("N43").Value = -50298

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1365.3334
$ws.Range("I96").Value = 1298
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 1298
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = 75
$ws.Range("N96").Value = -4246

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1775.5526
$ws.Range("I136").Value = 1334.6786
$ws.Range("J136").Value = 3010
$ws.Range("K136").Value = 4004.0358
$ws.Range("L136").Value = 9030
$ws.Range("M136").Value = -1454.0358
$ws.Range("N136").Value = -14130
